$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the (first) paragraph whose visible text, with all whitespace
# removed, equals the given normalized signature. Returns the Paragraph
# object, or $null if not found.
# ---------------------------------------------------------------------------
function Find-ParagraphByNormalizedText($normalizedTarget) {
    foreach ($p in $d.Paragraphs) {
        $cur = $p.Range.Text
        $curNorm = $cur -replace '\s', ''
        if ($curNorm -eq $normalizedTarget) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# Helper: find the (first) paragraph whose visible text contains the given
# substring (simple, case sensitive).
# ---------------------------------------------------------------------------
function Find-ParagraphContaining($needle) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like ("*" + $needle + "*")) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) The five "a-trou" oMath lines: collapse the 3-run whitespace-padded
#    equations into a single run, squashing the separate whitespace-only
#    runs away (keeping the same overall spacing/text).
# ---------------------------------------------------------------------------
function Set-OMathSingleRun($normalizedTarget, $newText) {
    foreach ($p in $d.Paragraphs) {
        $omaths = $p.Range.OMaths
        if ($omaths.Count -gt 0) {
            $om = $omaths.Item(1)
            $curNorm = $om.Range.Text -replace '\s', ''
            if ($curNorm -eq $normalizedTarget) {
                $xml = '<m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:lang w:eastAsia="fr-FR"/></w:rPr><m:t>' + $newText + '</m:t></m:r></m:oMath>'
                $om.Range.InsertXML($xml) | Out-Null
                return $true
            }
        }
    }
    return $false
}

Set-OMathSingleRun "0+=1"  "0+                 =1"  | Out-Null
Set-OMathSingleRun "0+=4"  "0+                 =4"  | Out-Null
Set-OMathSingleRun "0+=-2" "0+                 =-2" | Out-Null
Set-OMathSingleRun "2-=-1" "2-                 =-1" | Out-Null
Set-OMathSingleRun "2-=3"  "2-                 =3"  | Out-Null

# ---------------------------------------------------------------------------
# 2) "Nombres de " + "mêmes signes" (two runs) -> "Nombres de mêmes signes"
#    (one run), leaving the trailing " :" run untouched. The paragraph
#    keeps its "Titre2" style.
# ---------------------------------------------------------------------------
$p2 = Find-ParagraphContaining "Nombres de"
if ($p2 -ne $null -and $p2.Range.Text -like "*mêmes signes*") {
    $full = $p2.Range
    $xml = '<w:p><w:r><w:t>Nombres de mêmes signes</w:t></w:r><w:r><w:t xml:space="preserve"> :</w:t></w:r></w:p>'
    $full.InsertXML($xml) | Out-Null
    $p2b = Find-ParagraphByNormalizedText "Nombresdemêmessignes:"
    if ($p2b -ne $null) {
        $p2b.Style = "Titre2"
    }
}

# ---------------------------------------------------------------------------
# 3) "II" (title) before "Différence de deux nombres relatifs" gains a new
#    run "I" right after it, turning the displayed heading into "III".
# ---------------------------------------------------------------------------
$p3 = Find-ParagraphContaining "Différence de deux nombres relatifs"
if ($p3 -ne $null) {
    $full = $p3.Range
    $apos = [char]39
    $dash = [char]8211
    $xml = '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr></w:pPr>' +
           '<w:r><w:t>II</w:t></w:r>' +
           '<w:r><w:t>I</w:t></w:r>' +
           '<w:r><w:t xml:space="preserve"> ' + $dash + ' </w:t></w:r>' +
           '<w:r><w:t xml:space="preserve">Différence de deux nombres relatifs </w:t></w:r>' +
           '<w:r><w:t>:</w:t></w:r>' +
           '</w:p>'
    $full.InsertXML($xml) | Out-Null
    $p3b = Find-ParagraphContaining "Différence de deux nombres relatifs"
    if ($p3b -ne $null) {
        $p3b.Style = "Titre1"
    }
}

# ---------------------------------------------------------------------------
# 4) "III" (title, carries w:lastRenderedPageBreak) before "Simplification
#    d'une suite de sommes" splits into two runs "I" + "V", turning the
#    displayed heading into "IV".
# ---------------------------------------------------------------------------
$p4 = Find-ParagraphContaining "Simplification d"
if ($p4 -ne $null) {
    $full = $p4.Range
    $apos = [char]39
    $dash = [char]8211
    $xml = '<w:p><w:pPr><w:pStyle w:val="Titre1"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr></w:pPr>' +
           '<w:r><w:lastRenderedPageBreak/><w:t>I</w:t></w:r>' +
           '<w:r><w:t>V</w:t></w:r>' +
           '<w:r><w:t xml:space="preserve"> ' + $dash + ' </w:t></w:r>' +
           '<w:r><w:t>Simplification d' + $apos + 'une suite de sommes</w:t></w:r>' +
           '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
           '<w:r><w:t>:</w:t></w:r>' +
           '</w:p>'
    $full.InsertXML($xml) | Out-Null
    $p4b = Find-ParagraphContaining "Simplification d"
    if ($p4b -ne $null) {
        $p4b.Style = "Titre1"
    }
}
